$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last data row (row 27, "Vnzlno Punta Anca" duplicate-shifted entry /
# originally "Santander Norte" row) so the table shrinks from 26 to 25 data rows.
$ws.Rows(27).Delete()

# Updated client names (normalized to upper case, trimmed) together with the
# refreshed Fecha (serial date) and Valor figures for each remaining row.
$data = @(
  @{Row=2;  Name="ALISO";                   Fecha=45983; Valor=101000},
  @{Row=3;  Name="CAMPO VERDE ZIPAQUIRA";    Fecha=45988; Valor=64200},
  @{Row=4;  Name="CARNES JOHANA";            Fecha=45993; Valor=176000},
  @{Row=5;  Name="CARNILANDIA";              Fecha=45994; Valor=436700},
  @{Row=6;  Name="CARNIVOROS";               Fecha=45959; Valor=200000},
  @{Row=7;  Name="CIMARRON DORADO";          Fecha=45992; Valor=407000},
  @{Row=8;  Name="COCINA CHINA";             Fecha=45992; Valor=170000},
  @{Row=9;  Name="COCINA CHINA";             Fecha=45994; Valor=85000},
  @{Row=10; Name="DARWIN FUTBOL";            Fecha=45921; Valor=200000},
  @{Row=11; Name="DAVIDCITO";                Fecha=45947; Valor=100000},
  @{Row=12; Name="EL RUBY";                  Fecha=45992; Valor=85100},
  @{Row=13; Name="LA PAMPA";                 Fecha=45994; Valor=249000},
  @{Row=14; Name="LA SELECTA";               Fecha=45912; Valor=82000},
  @{Row=15; Name="MARIANA";                  Fecha=45650; Valor=171900},
  @{Row=16; Name="MERKA FRUVER ALEJANDRO";   Fecha=45988; Valor=60900},
  @{Row=17; Name="MERKA FRUVER MILDRED";     Fecha=45988; Valor=115400},
  @{Row=18; Name="MEZA 2";                   Fecha=45989; Valor=188000},
  @{Row=19; Name="MULTICARNES";              Fecha=45989; Valor=558300},
  @{Row=20; Name="NOVILLON SAN MATEO";       Fecha=45971; Valor=83000},
  @{Row=21; Name="PINILLA";                  Fecha=45931; Valor=166000},
  @{Row=22; Name="PINILLA";                  Fecha=45924; Valor=16000},
  @{Row=23; Name="PINILLA SOACHA";           Fecha=45993; Valor=129000},
  @{Row=24; Name="PLAZA JESSICA";            Fecha=45993; Valor=621000},
  @{Row=25; Name="SANTANDER SUR";            Fecha=45993; Valor=80000},
  @{Row=26; Name="VNZLNO PUNTA ANCA";        Fecha=45992; Valor=82000}
)

foreach ($item in $data) {
  $r = $item.Row
  $ws.Cells.Item($r, 2).Value = $item.Name
  $ws.Cells.Item($r, 3).Value = $item.Fecha
  $ws.Cells.Item($r, 4).Value = $item.Valor
}
